# Apply updated input data (columns B:G, rows 3-11) and scroll position change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New input values for B:G across rows 3-11 (downstream formula columns
# P:AJ recompute automatically since they are formulas referencing these cells).
$data = @{
    3  = @(107, 206, 70, 232, 0, 124)
    4  = @(118, 218, 66, 227, 23, 144)
    5  = @(129, 300, 101, 415, 15, 149)
    6  = @(87, 152, 40, 198, 0, 90)
    7  = @(83, 166, 52, 210, 0, 116)
    8  = @(145, 214, 105, 268, 0, 112)
    9  = @(130, 240, 87, 251, 7, 134)
    10 = @(111, 217, 72, 235, 0, 120)
    11 = @(103, 232, 80, 240, 8, 132)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $values[$i]
    }
}

# Restore the previously scrolled viewport (top-left visible cell).
$ws.Application.ActiveWindow.ScrollColumn = 26
$ws.Application.ActiveWindow.ScrollRow = 1

$wb.Save()
